# Update the build timestamp embedded in the "version" strings from
# "February 03 2026 17.29.55 EST" to "February 03 2026 18.05.36 EST"
# across the "About" sheet and the "Boundaries and methane sources" sheet.

$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")

# A2 and A6 on the "About" sheet contain the version/citation strings.
$aboutCells = @("A2", "A6")
foreach ($addr in $aboutCells) {
    $cell = $aboutSheet.Range($addr)
    $text = $cell.Value()
    if ($text -ne $null -and $text.Contains($oldStamp)) {
        $cell.Value = $text.Replace($oldStamp, $newStamp)
    }
}

# S2:S53 on the data sheet hold the per-row build_version string.
for ($row = 2; $row -le 53; $row++) {
    $cell = $dataSheet.Cells.Item($row, 19)  # column S = 19
    $text = $cell.Value()
    if ($text -ne $null -and $text.Contains($oldStamp)) {
        $cell.Value = $text.Replace($oldStamp, $newStamp)
    }
}
